$d = $word.ActiveDocument

# Locate the paragraph ending in "...started what would come to be known as the Browser Wars."
# (the last real content paragraph before the trailing empty ones / sectPr).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*started what would come to be known as the Browser Wars.*") {
        $target = $cand
    }
}

if ($null -eq $target) {
    throw "Could not locate target paragraph"
}

$start = $target.Range.Start
$end = $target.Range.End - 1
$full = $d.Range($start, $end)

$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Netscape, however, would ultimately fail to continue the market dominance it had enjoyed. Microsoft, </w:t></w:r><w:r><w:t>and its CEO at the time, Bill Gates, was worried that platform-independent browsers like Netscape threatened the dominance of their flagship Windows operating system and were determined to take on Andreessen and his team. Licensing the browser created by Spyglass, whose own software was largely composed of Andreessen’s original Mosaic browser, the Windows team built Internet Explorer, and in doing so, started what would come to be known as the Browser Wars.</w:t></w:r><w:r><w:t xml:space="preserve"> These Browser Wars would push Netscape, and Andreessen out of business after a long struggle of features, dirty dealing by Microsoft</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> and a litany of lawsuits that nearly saw Microsoft broken up through anti-trust legislation.</w:t></w:r><w:r><w:t xml:space="preserve"> The decline of Netscape during the Browser Wars led to its acquisition by AOL in 1999 for more than </w:t></w:r><w:r><w:t>4.3 billion dollars</w:t></w:r><w:r><w:t>, an immense sum especially for the time, and Marc Andreessen departed Netscape.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:tab/><w:t xml:space="preserve">When Andreessen left, he did not end his software career and his impact on the industry at large was far from it. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Again</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Andreessen had an innovative and impactful idea that changed the whole of the software industry, cloud computing. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Opsware</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, originally called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Loudcloud</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, was the result of Andreesen and fellow software innovator Ben Horowitz’s belief that </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>computing</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and software services would be the </w:t></w:r><w:r><w:t xml:space="preserve">future of </w:t></w:r><w:r><w:t>consumer-facing</w:t></w:r><w:r><w:t xml:space="preserve"> e-commerce companies. </w:t></w:r><w:r><w:t>Andreessen</w:t></w:r><w:r><w:t xml:space="preserve"> pioneered the idea of software as a service (SaaS) for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Opsware’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> cloud computing systems</w:t></w:r><w:r><w:t xml:space="preserve">, and by June of 2000, raised over $120 million during its second round of funding, the largest amount ever raised in a Series B at that time. </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:tab/></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Andreessen’s work at </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Opsware</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> would lead to many of the features that have let cloud computing become what it is today. Notably</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> they were one of the first companies to integrate their datacenters and features end-to-end so that customers could handle </w:t></w:r><w:r><w:t>all</w:t></w:r><w:r><w:t xml:space="preserve"> their cloud computing needs under the same platform. Andreessen accomplished this through a series of acquisitions and integrations </w:t></w:r><w:r><w:t>that allowed them to eventually create enough value to be acquired by Hewlett-Packard in 2007 for over $1.65 billion in cash, almost sixteen times its average annual revenue, and formed a large part of Hewlett-Packard Enterprise, HP’s flagship B2B SaaS offering.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Again, this latest company exit was not the end of Andreessen’s impact on the software industry. This time, rather than founding companies, Andreessen and his fellow </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Opsware</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Co-founder, Ben Horowitz, would create Andreessen Horowitz, one of the most notable and successful Silicon Valley venture capital firms in history. Starting with a capitalization of $300 million, the pair quickly grew the firm’s assets under management to more than $1.2 billion in under 2 years, by which time they had invested in a </w:t></w:r><w:r><w:t xml:space="preserve">series of now notable companies including Okta, Skype, Facebook (now Meta), Groupon, Twitter, Zynga, Airbnb, and Stripe. </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:tab/></w:r></w:p><w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">During 2012-2013, Andreessen Horowitz directed investments in some of the most notable, </w:t></w:r><w:r><w:t>up-and-coming</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>startups</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">in a variety of technology niches, many of which are only now coming to the fore. </w:t></w:r><w:r><w:t xml:space="preserve">Some of the most notable investments made at this time are in the cryptocurrency space, </w:t></w:r><w:r><w:t>and</w:t></w:r><w:r><w:t xml:space="preserve"> are now proving to have been visionary investments, especially that of Coinbase, the largest cryptocurrency platform in the US, and Ripple a blockchain payments platform that created one of the largest cryptocurrencies by market capitalization, XRP. Additionally, Andreessen and his firm invested around $100 million into GitHub and made significant investments in Lyft and Oculus during this time.  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$full.InsertXML($frag)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
